# Applies the scheduled-runner profit recalculation update across all leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 557.38464
$ws.Range("I33").Value = 194.6
$ws.Range("J33").Value = 1766.6666
$ws.Range("K33").Value = 194.6
$ws.Range("L33").Value = 1766.6666
$ws.Range("M33").Value = 34.40000000000001
$ws.Range("N33").Value = -2224.6666
$ws.Range("H41").Value = 1300.6666
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 1902
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 1902
$ws.Range("M41").Value = -560
$ws.Range("N41").Value = -2782
$ws.Range("H51").Value = 5297.6665
$ws.Range("J51").Value = 5297.6665
$ws.Range("L51").Value = 5297.6665
$ws.Range("N51").Value = -6265.6665
$ws.Range("H74").Value = 4608.3335
$ws.Range("I74").Value = 2450
$ws.Range("J74").Value = 5040
$ws.Range("K74").Value = 2450
$ws.Range("L74").Value = 5040
$ws.Range("M74").Value = -1514
$ws.Range("N74").Value = -6912
$ws.Range("H77").Value = 4608.3335
$ws.Range("I77").Value = 2450
$ws.Range("J77").Value = 5040
$ws.Range("K77").Value = 12250
$ws.Range("L77").Value = 25200
$ws.Range("M77").Value = -7570
$ws.Range("N77").Value = -34560
$ws.Range("H107").Value = 761.75757
$ws.Range("I107").Value = 534.0417
$ws.Range("K107").Value = 534.0417
$ws.Range("M107").Value = 1385.9583
$ws.Range("H112").Value = 2854.3635
$ws.Range("I112").Value = 1132.6666
$ws.Range("K112").Value = 3397.9998
$ws.Range("M112").Value = -2289.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1492.3478
$ws.Range("I110").Value = 684.6667
$ws.Range("J110").Value = 4400
$ws.Range("K110").Value = 684.6667
$ws.Range("L110").Value = 4400
$ws.Range("M110").Value = 1360.3333
$ws.Range("N110").Value = -8490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H86").Value = 2568.6667
$ws.Range("I86").Value = 2568.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2568.6667
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1445.6667
$ws.Range("H89").Value = 2568.6667
$ws.Range("I89").Value = 2568.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12843.3335
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -7227.333500000001
$ws.Range("H107").Value = 812.8461
$ws.Range("I107").Value = 774.2727
$ws.Range("J107").Value = 1025
$ws.Range("K107").Value = 774.2727
$ws.Range("L107").Value = 1025
$ws.Range("M107").Value = 1145.7273
$ws.Range("N107").Value = -4865
$ws.Range("N9").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2685.3076
$ws.Range("I31").Value = 1403.3334
$ws.Range("K31").Value = 1403.3334
$ws.Range("M31").Value = -1108.3334
$ws.Range("H34").Value = 2685.3076
$ws.Range("I34").Value = 1403.3334
$ws.Range("K34").Value = 1403.3334
$ws.Range("M34").Value = -1201.3334
$ws.Range("H99").Value = 2083.1667
$ws.Range("I99").Value = 2575
$ws.Range("J99").Value = 1099.5
$ws.Range("K99").Value = 2575
$ws.Range("L99").Value = 1099.5
$ws.Range("M99").Value = -1077
$ws.Range("N99").Value = -4095.5
$ws.Range("H126").Value = 2083.1667
$ws.Range("I126").Value = 2575
$ws.Range("J126").Value = 1099.5
$ws.Range("K126").Value = 7725
$ws.Range("L126").Value = 3298.5
$ws.Range("M126").Value = -5255
$ws.Range("N126").Value = -8238.5
$ws.Range("H134").Value = 4741.4287
$ws.Range("I134").Value = 4417.4
$ws.Range("K134").Value = 13252.2
$ws.Range("M134").Value = -10717.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 401999.2
$ws.Range("I11").Value = 2000000
$ws.Range("K11").Value = 6000000
$ws.Range("M11").Value = -5999860
$ws.Range("H86").Value = 378.85715
$ws.Range("I86").Value = 297.5
$ws.Range("J86").Value = 411.4
$ws.Range("K86").Value = 892.5
$ws.Range("L86").Value = 1234.2
$ws.Range("M86").Value = 293.5
$ws.Range("N86").Value = -3606.2
$ws.Range("H89").Value = 378.85715
$ws.Range("I89").Value = 297.5
$ws.Range("J89").Value = 411.4
$ws.Range("K89").Value = 2677.5
$ws.Range("L89").Value = 3702.6
$ws.Range("M89").Value = 3250.5
$ws.Range("N89").Value = -15558.6
$ws.Range("H137").Value = 7498.6665
$ws.Range("J137").Value = 8598.4
$ws.Range("L137").Value = 25795.2
$ws.Range("N137").Value = -35995.2
$ws.Range("H138").Value = 1613.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3007.5386
$ws.Range("J16").Value = 1933.3334
$ws.Range("L16").Value = 1933.3334
$ws.Range("N16").Value = -2273.3334
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H93").Value = 797
$ws.Range("I93").Value = 797
$ws.Range("K93").Value = 797
$ws.Range("M93").Value = 451
$ws.Range("H106").Value = 11272
$ws.Range("J106").Value = 11272
$ws.Range("L106").Value = 11272
$ws.Range("N106").Value = -13796
$ws.Range("H132").Value = 17022.75
$ws.Range("I132").Value = 21796.166
$ws.Range("J132").Value = 2702.5
$ws.Range("K132").Value = 65388.49800000001
$ws.Range("L132").Value = 8107.5
$ws.Range("M132").Value = -62858.49800000001
$ws.Range("N132").Value = -13167.5
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1172.6666
$ws.Range("I96").Value = 960
$ws.Range("J96").Value = 1598
$ws.Range("K96").Value = 960
$ws.Range("L96").Value = 1598
$ws.Range("M96").Value = 413
$ws.Range("N96").Value = -4344
$ws.Range("H113").Value = 2683.3
$ws.Range("I113").Value = 1262.5714
$ws.Range("K113").Value = 3787.7142
$ws.Range("M113").Value = -1617.7142
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("H135").Value = 45499.875
$ws.Range("I135").Value = 19499.5
$ws.Range("J135").Value = 54166.668
$ws.Range("K135").Value = 19499.5
$ws.Range("L135").Value = 54166.668
$ws.Range("M135").Value = -14429.5
$ws.Range("N135").Value = -64306.668
